$d = $word.ActiveDocument

function Set-ParaText($para, [string]$text) {
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

# 1) Straightforward text replacements that keep paragraph structure/style as-is
Set-ParaText $d.Paragraphs(1) '數學 - 盈虧問題 - 基本型'
Set-ParaText $d.Paragraphs(2) '(1) 貓媽媽給小貓分魚，每隻小貓分10條魚，就多出8條魚，每隻小貓分11條魚則正好分完。問：一共有【8】隻小貓、貓媽媽一共有【88】條魚'
Set-ParaText $d.Paragraphs(3) '(2) 童軍團去種樹，如果每人種6棵，則少13棵樹苗；如果每人種5棵，則正好種完。問：童軍團有【13】個人，有【65】棵樹苗'
Set-ParaText $d.Paragraphs(4) '(3) 臥龍自然保護區管理員把一些竹子分給若干隻大熊貓，每隻大熊貓分5個還多餘10棵竹子，如果大熊貓數增加到3倍還少5隻，那麼每隻大熊貓分2棵竹子還缺少8棵竹子，問有大熊貓【28】隻，竹子【150】棵'
Set-ParaText $d.Paragraphs(5) '(4) 小紅家買來一籃桔子，分給全家人。如果其中二人每人分4個，其餘每人分2個，那麼多出4個；如果一人分6個，其餘每人分4個，那麼缺12個。問：小紅家買來【26】個桔子、小紅家共有【 9】人'
Set-ParaText $d.Paragraphs(6) '(5) 一些桔子分給若干人，每人5個餘10個桔子，如果人數增加到3倍還少5人，那麼每人分2個還缺8個，有桔子【150】個。'
Set-ParaText $d.Paragraphs(7) '(6) 學生栽植一批樹苗，如果每個學生栽6棵，剩12棵；如果其中9個學生每人栽4棵，而其餘學生栽8棵，結果缺2棵。問：這批樹苗有【162】棵、參加植樹的學生有【25】人'
Set-ParaText $d.Paragraphs(8) '(7) 小朋友分蘋果，每人分18個，還多出2個；每人分20個，就有一位小朋友沒分到蘋果，問：共有【11】個小朋友、共有【200】個蘋果'
Set-ParaText $d.Paragraphs(9) '(8) 48本書分給兩組小朋友，已知第二組比第一組多5人。如果把書全部分給第一組，那麼每人4本，有剩餘；每人5本，書不夠。如果把書全分給第二組，那麼每人3本，有剩餘；每人4本，書不夠，問：第一組有【10】人、第二組有【15】人'
Set-ParaText $d.Paragraphs(10) '(9) 老師說：YEARDFASF有【214】人JJJ有【999】個'
Set-ParaText $d.Paragraphs(11) '(10) 幼稚園將一筐蘋果分給大班和小班的小朋友，如果大班每人分5個，就多10個；如果小班每人分8個，就少了2個。已知大班比小班多3人。問：這筐蘋果有【70】個'
Set-ParaText $d.Paragraphs(12) '(11) 同學們暑假前到圖書館借書，如果每人借4本，則最後少2本；如果前2人每人先借8本，餘下的人每人借3本，這些圖書恰好借完。問：書的總數是【46】本'
Set-ParaText $d.Paragraphs(13) '(12) 媽媽買回一筐蘋果，按計劃吃的天數算了一下，如果每天吃4個，要多出48個蘋果；如果每天吃6個，則多出8個蘋果。問：媽媽買回的蘋果有【 128 】個，計畫吃【 20 】天'
Set-ParaText $d.Paragraphs(14) '(13) 三年級一班少先隊員參加學校搬磚勞動，如果每人搬4塊磚，還剩7塊；如果每人搬5塊，則多1塊磚。問：這個班少先隊有【 6 】個人，磚共有【 31 】塊'
Set-ParaText $d.Paragraphs(15) '(14) 媽媽買回一筐蘋果，按計劃吃的天數算了一下，如果每天吃4個，要多出48個蘋果；如果每天吃6個，則又少8個蘋果。問：媽媽買回的蘋果有【160】個、計畫吃【28】天 '
Set-ParaText $d.Paragraphs(18) '(16) 學校買來一批電風扇分給各班。若有兩個班每班分到4台，其餘每班只能分2台；如果有一個班分6台，其餘每班分4台，還差12台。問：共買來【18】台電風扇、有【7】個班'
Set-ParaText $d.Paragraphs(19) '(17) 將蘋果放入一些籃子中，如果每籃放8個，則缺少21個；如果每籃改為放6個，則缺少3個。問：籃子有【9】個，蘋果有【51】個'

# 2) New question (15) replaces the old (15) text paragraph + the picture paragraph.
#    Insert a fresh paragraph right after (14) so it inherits the 'question' style / left jc,
#    then delete the two old paragraphs (old text (15) and the picture paragraph).
$p14 = $d.Paragraphs(15)
$p14.Range.InsertParagraphAfter()
Set-ParaText $d.Paragraphs(16) '(15) 明明過生日，同學們去給他買蛋糕，如果每人出8元，就多出8元；每人出7元，就多出4元。問：有【  4  】個同學去買蛋糕、這個蛋糕的價錢【  24  】元'
$d.Paragraphs(17).Range.Delete()
$d.Paragraphs(17).Range.Delete()

# 3) Append the 23 new questions (18)-(40) at the end of the document.
$newQuestions = @(
    '(18) 老師把一袋糖分給小朋友。如果只分給小班，每人可得12塊，如果分給中班和小班，每人只能分到4塊。如果這袋糖只分給中班，每人可分到【 6 】塊',
    '(19) 小芳把鮮花插入一些花瓶中，如果每個花瓶裡插5枝，則多12枝；如果每個花瓶裡插8枝，還多3枝。問：每個花瓶裡插【9】枝花可以剛好把鮮花分完',
    '(20) 某廠運來一批煤，如果每天燒1500公斤，那麼比原計劃提前一天燒完；如果每天燒1000公斤，那麼將比原計劃多用一天。問：現在要求按原計劃燒完，每天應燒煤【1200】公斤',
    '(21) 一個商販估計，若1公斤蘋果賣24元，會賠40元；若1公斤蘋果賣30元，可以賺80元。問：若以不賠不賺的價格賣出，每公斤蘋果應賣【26】元',
    '(22) 幼稚園老師給小朋友分梨子，如果每人分4個，則多9個；如果每人分5個，則少6個。問：有【15】個小朋友，有【69】個梨子',
    '(23) 老師把一籃蘋果分給小朋友，如果減少一名同學，每個同學正好分得5個；如果增加一個同學，正好每人分得4個，求這籃蘋果一共有【 40 】個',
    '(24) 工人鋪一條路基，若每天鋪260公尺，鋪完全路長就得延長8天；若每天鋪300公尺，鋪完全路長仍要延長4天。問：這條路長【7800】公尺。',
    '(25) 老師把一些蘋果分給小朋友。如果每人分一個，還剩下8個蘋果；如果每人分2個，那麼還少2個蘋果。問：共有【10】個小朋友',
    '(26) 甲、乙兩組同學做紅花，每人做8朵，正好送給五年級每個同學一朵。如果把這些紅花讓甲組單獨做，每人要多做4朵。如果把這些紅花讓乙組同學單獨做，每人要做【24】朵',
    '(27) 四年級某班的同學去植樹，他們分了一下小組，如果增加一個小組，正好每小組5人；如果減少一小組，正好每組7人。問：這個班共有【35】人',
    '(28) 小朋友分糖果，每人3粒，餘30粒；每人5粒，少4粒。問：有【17】個小朋友、【81】粒糖',
    '(29) 302班同學去植樹，每人植5棵，有3棵沒有人植，如果其中4人每人植4棵，其餘每人植6棵，就恰好植完所有的樹。問：參加植樹的有【11】名同學、共植【58】棵樹',
    '(30) 幼稚園的小朋友分玩具，如果每人分18件，還多出2件；如果每人分20件，就有1位小朋友分不到玩具。問：共有【11】個小朋友、【200】件玩具',
    '(31) 體育隊將一些羽毛球分給若干個人，每人5個還多餘10個羽毛球，如果人數增加到 3倍，那麼每人分2個羽毛球還缺少8個，問：有羽毛球【100】個',
    '(32) 小朋友分糖果，每人分10粒，正好分完；若每人分16粒，則有3個小朋友分不到糖果。問：有【80】粒糖果',
    '(33) 食堂採購員小李去買肉，如果買牛肉18公斤，那麼差40元；如果買豬肉20公斤，那麼多20元。已知牛肉比豬肉每公斤貴8元。問：牛肉每公斤【50】元、豬肉每公斤【42】元，小李帶了【860】元',
    '(34) 農民鋤草，其中5人各鋤4畝，餘下的各鋤3畝，這樣分配最後餘下26畝；如果其中3人每人各鋤3畝，餘下的人各鋤5畝，最後餘下3畝。問：草地面積【82】畝、鋤草人數【17】人',
    '(35) 三年級給優秀學生發獎品書，如果每個學生發5冊還剩32冊；如果其中10個學生每人發4冊，其餘每人發8冊，就恰好發完。問：優秀學生有【24】人、獎品書有【152】冊',
    '(36) 老師將一些練習本發給班上的學生。如果每人發10本，則有兩個學生沒分到；如果每人發8本，則正好發完。問：有【10】個學生、【80】本練習本',
    '(37) 學校園林科有一批樹苗，交給若干名學生去栽，一次一次往下分，每次分一棵，最後剩下12棵，不夠分了。如果再拿來8棵，那麼每個學生正好栽10棵。問：參加栽樹的學生有【20】人，這批樹苗共【192】棵',
    '(38) 老師把一籃蘋果分給小朋友，如果減少一名同學，每個同學正好分得5個；如果增加一個同學，正好每人分得4個。問：小朋友【9】人，蘋果共有【45】個',
    '(39) 猴王帶領一群猴子去摘桃。下午收工後，猴王開始分配，若大猴分5個，小猴分3個，猴王可留10個；若大、小猴都分4個，猴王能留下20個。在這群猴子中，大猴（不包括猴王）比小猴多【  10 】隻。',
    '(40) 幼稚園小朋友分蘋果，如果每人分3個就多了11個，如果每人分5個還缺5個。問：有【8】個小朋友、蘋果有【35】個'
)
$last = $d.Paragraphs($d.Paragraphs.Count)
foreach ($q in $newQuestions) {
    $last.Range.InsertParagraphAfter()
    $last = $d.Paragraphs($d.Paragraphs.Count)
    Set-ParaText $last $q
}

Write-Host ("Final paragraph count: " + $d.Paragraphs.Count)
